$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 ("I0") and J1 ("IF") ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell H1 (bold, thin border,
# centered/top aligned) onto the two new header cells so they match the
# look of the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows I2:J23 ---
$data = @(
    @(5,5),
    @(8,8),
    @(6,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(5,5),
    @(9,9),
    @(6,6),
    @(5,7),
    @(5,6),
    @(9,9),
    @(4,7),
    @(3,4),
    @(7,8),
    @(6,6),
    @(8,8),
    @(4,4),
    @(6,6),
    @(4,4),
    @(5,5),
    @(2,2)
)

$r = 2
foreach ($pair in $data) {
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
    $r = $r + 1
}
